# Add the missing "John Cena" sample record across the three sheets
# (students, class_data, mcas_scores) to support the new admin page.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# students sheet: add new student record (John Cena) in row 3
# ---------------------------------------------------------------------------
$wsStudents = $wb.Worksheets.Item("students")

$wsStudents.Range("A3").Value = 125687
$wsStudents.Range("C3").Value = "John"
$wsStudents.Range("B3").Value = "Cena"
$wsStudents.Range("D3").Value = "Wrestling"
$wsStudents.Range("M3").Value = "MA"
$wsStudents.Range("N3").Value = "USA"
$wsStudents.Range("O3").Value = 44338
$wsStudents.Range("P3").Value = "Graduated"
$wsStudents.Range("Q3").Value = 158
$wsStudents.Range("R3").Value = 2.9
$wsStudents.Range("S3").Value = 2.99
$wsStudents.Range("T3").Value = 3.54
$wsStudents.Range("U3").Value = 800
$wsStudents.Range("V3").Value = 21
$wsStudents.Range("W3").Value = "Y"
$wsStudents.Range("X3").Value = "White"

# ---------------------------------------------------------------------------
# mcas_scores sheet: add John Cena's MCAS score row
# ---------------------------------------------------------------------------
$wsMcas = $wb.Worksheets.Item("mcas_scores")

$wsMcas.Range("A3").Value = 125687
$wsMcas.Range("B3").Value = 22
$wsMcas.Range("C3").Value = 198
$wsMcas.Range("D3").Value = "NI"

$null = $wsMcas.Range("F13").Select()
$null = $wsStudents.Range("B39").Select()

# ---------------------------------------------------------------------------
# class_data sheet: add John Cena's wrestling course history (rows 5-8)
# ---------------------------------------------------------------------------
$wsClass = $wb.Worksheets.Item("class_data")

$wsClass.Range("A5").Value = 125687
$wsClass.Range("B5").Value = "UNDG"
$wsClass.Range("C5").Value = "Day - Science & Eng "
$wsClass.Range("D5").Value = "Intro to Wrestling"
$wsClass.Range("G5").Value = "FA"
$wsClass.Range("H5").Value = 2019

$wsClass.Range("A6").Value = 125687
$wsClass.Range("B6").Value = "UNDG"
$wsClass.Range("C6").Value = "Day - Science & Eng "
$wsClass.Range("D6").Value = "Wrestling 2"
$wsClass.Range("G6").Value = "SP"
$wsClass.Range("H6").Value = 2020

$wsClass.Range("E5").Value = "WES1001"
$wsClass.Range("E6").Value = "WES1002"

$wsClass.Range("A7").Value = 125687
$wsClass.Range("B7").Value = "UNDG"
$wsClass.Range("C7").Value = "Day - Science & Eng "
$wsClass.Range("D7").Value = "Top Rope Jumping"
$wsClass.Range("G7").Value = "FA"
$wsClass.Range("H7").Value = 2021

$wsClass.Range("E7").Value = "WES4955"

$wsClass.Range("A8").Value = 125687
$wsClass.Range("B8").Value = "UNDG"
$wsClass.Range("C8").Value = "Day - Science & Eng "
$wsClass.Range("D8").Value = "Wresting Capstone"
$wsClass.Range("G8").Value = "SP"
$wsClass.Range("H8").Value = 2021

$wsClass.Range("E8").Value = "WES9999"

$wsClass.Range("F5").Value = "A"
$wsClass.Range("F8").Value = "F"
$wsClass.Range("F7").Value = "B"
$wsClass.Range("F6").Value = "C-"

# Leave class_data as the active sheet/selection, matching the saved state.
$null = $wsClass.Range("F5").Select()
